$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from an existing
# header cell (A1) onto the three new header cells so they match the
# rest of row 1's formatting (style index 1 in the original file).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels for the season-record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate every data row (2-55) with the team's season record.
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value = 53
    $ws.Cells.Item($r, 31).Value = 109
    $ws.Cells.Item($r, 32).Value = 0
}
